# "Removed indexing, since it broke the code if there was index and no
# name" — drop the running-index values that lived in column A (rows 2
# through 142; the header row and the trailing rows that never had an
# index are left alone). Columns B (name) and C (wishes) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A142").ClearContents()

# Leave the cursor parked near the top of the now-reindexed sheet.
$ws.Range("B4").Select()
